$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (losing trailing zeros / exact text form).
$textPriceCells = @("D5", "D6", "D12", "D13", "D14", "D17", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D41", "D46", "D48", "D49", "D51")
foreach ($cellAddr in $textPriceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "72.116.47"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "4.036.37"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "540.38"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "151.66"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").Value = "4.032.01"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "53.86"
$ws.Range("E12").Value = "  +10.86%  "
$ws.Range("D13").Value = "0.0000328"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "10.89"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "4.682.94"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "4.044.27"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "14.38"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "72.116.44"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "447.31"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("D23").Value = "97.51"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").Value = "14.64"
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "4.32"
$ws.Range("E27").Value = "  +17.47%  "
$ws.Range("D28").Value = "11.32"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "10.79"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").Value = "5.96"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").Value = "37.14"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  +18.31%  "
$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "13.60"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "49.18"
$ws.Range("E35").Value = "  +14.35%  "
$ws.Range("D36").Value = "679.21"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "66.82"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "0.457"
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("D39").Value = "0.0₃0875"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("E43").Value = "  +16.98%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "3.11"
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "3.43"
$ws.Range("E51").Value = "  +1.27%  "
